# Add Revisi Nota (UNTESTED)
# Adds four new TODO entries to the "TODO" worksheet, including one
# entry ("Edit Sales Order, memanfaatkan copy nota") that is bold and
# has a bold space character between "," and "memanfaatkan".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")

# B9: plain text entry
$ws.Range("B9").Value = "Simpan data orang yg mengedit sales order"

# B10: plain text entry
$ws.Range("B10").Value = "Copy Nota menampilkan data orang yg mengedit sales order"

# B8: bold cell style, with just the space between "," and "memanfaatkan" bold
$ws.Range("B8").Value = "Edit Sales Order, memanfaatkan copy nota"
$ws.Range("B8").Font.Bold = $true
$ws.Range("B8").Characters(17, 1).Font.Bold = $false
$ws.Range("B8").Characters(18, 1).Font.Bold = $true
$ws.Range("B8").Characters(19, 22).Font.Bold = $false

# B13: plain text entry (rows 11-12 stay empty)
$ws.Range("B13").Value = "Tambah persenan untuk harga di master produk"

# Move / record the active selection as it would be after typing the last entry
$ws.Range("B14").Select() | Out-Null

# Match page setup change captured in the diff
$ws.PageSetup.Orientation = 1
